# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (column D) and
# "Correspond Handback DateTime" (column G) for the first data row
# (the 851bafc6-... record) on both the "zh-cn" and "de-de" sheets,
# reflecting a newly generated handback report run.

$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# zh-cn sheet, row 2 (851bafc6-... record)
$wsZh.Range("D2").Value = "2016-02-17 09:44:03"
$wsZh.Range("G2").Value = "2016-02-17 09:44:49"

# de-de sheet, row 2 (851bafc6-... record)
$wsDe.Range("D2").Value = "2016-02-17 09:44:19"
$wsDe.Range("G2").Value = "2016-02-17 09:45:12"
